$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("B2").Value = "BRUNO125"
$ws.Range("B3").Value = "Marcela110"
$ws.Range("B4").Value = "BrUn97"
$ws.Range("B5").Value = "Maria25"

$ws.Range("B5").Select()
